$wb = $excel.ActiveWorkbook

# ALC row 33: Glazed and Confused | Clear Glass Lens (item 5512)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 321.4
$ws.Range("I33").Value = 246
$ws.Range("K33").Value = 246
$ws.Range("M33").Value = -17

# ALC row 99: Rumor Has It | Commanding Craftsman's Tea (item 19883)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H99").Value = 413.66666
$ws.Range("I99").Value = 256.4
$ws.Range("J99").Value = 1200
$ws.Range("K99").Value = 769.1999999999999
$ws.Range("L99").Value = 3600
$ws.Range("M99").Value = 728.8000000000001
$ws.Range("N99").Value = -6596

# ARM row 32: Ingot We Trust | Steel Ingot (item 44147)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 8141.8594
$ws.Range("I32").Value = 6041.1523
$ws.Range("K32").Value = 6041.1523
$ws.Range("M32").Value = -5754.1523

# ARM row 45: Hollow Hallmarks | Mythril Ingot (item 27714)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 1100
$ws.Range("I45").Value = 1166.6666
$ws.Range("K45").Value = 1166.6666
$ws.Range("M45").Value = -789.6666

# ARM row 74: As the Bolt Flies | Titanium Nugget (item 44000)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 3582.6
$ws.Range("I74").Value = 615.34784
$ws.Range("K74").Value = 615.34784
$ws.Range("M74").Value = 258.65216

# ARM row 77: Heavy Metal Banned (L) | Titanium Nugget (item 44000)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 3582.6
$ws.Range("I77").Value = 615.34784
$ws.Range("K77").Value = 3076.7392
$ws.Range("M77").Value = 1291.2608

# BSM row 20: Smelt and Dealt | Iron Ingot (item 14149)
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 2399.1177
$ws.Range("I20").Value = 2270.3572
$ws.Range("K20").Value = 2270.3572
$ws.Range("M20").Value = -2023.3572

# BSM row 43: Don't Fear the Reaper | Steel Scythe (item 22904)
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H43").Value = 200000
$ws.Range("J43").Value = 200000
$ws.Range("L43").Value = 200000
$ws.Range("N43").Value = -200362

# BSM row 56: I'd Rather Be Digging | Electrum Head Knife (item 2427)
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H56").Value = 40666.832
$ws.Range("J56").Value = 50000.25
$ws.Range("L56").Value = 50000.25
$ws.Range("N56").Value = -51478.25

# BSM row 80: Unbreaker | Titanium Ingot (item 13747)
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H80").Value = 272.57144
$ws.Range("I80").Value = 356.33334
$ws.Range("J80").Value = 249.72728
$ws.Range("K80").Value = 356.33334
$ws.Range("L80").Value = 249.72728
$ws.Range("M80").Value = 641.66666
$ws.Range("N80").Value = -2245.72728

# BSM row 83: Attack on Titanium (L) | Titanium Ingot (item 13747)
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H83").Value = 272.57144
$ws.Range("I83").Value = 356.33334
$ws.Range("J83").Value = 249.72728
$ws.Range("K83").Value = 1781.6667
$ws.Range("L83").Value = 1248.6364
$ws.Range("M83").Value = 3210.3333
$ws.Range("N83").Value = -11232.6364

# BSM row 92: Have Blade, Will Travel | High Steel Katzbalger (item 18033)
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H92").Value = 23075.8
$ws.Range("J92").Value = 23075.8
$ws.Range("L92").Value = 23075.8
$ws.Range("N92").Value = -28067.8

# BSM row 107: The Gold Experience | Deepgold Nugget (item 27706)
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 1600.625
$ws.Range("I107").Value = 1644.2858
$ws.Range("J107").Value = 1566.6666
$ws.Range("K107").Value = 1644.2858
$ws.Range("L107").Value = 1566.6666
$ws.Range("M107").Value = 275.7141999999999
$ws.Range("N107").Value = -5406.6666

# CRP row 31: Wall Not Found | Walnut Lumber (item 44023)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1514
$ws.Range("I31").Value = 1482.2858
$ws.Range("J31").Value = 1602.8
$ws.Range("K31").Value = 1482.2858
$ws.Range("L31").Value = 1602.8
$ws.Range("M31").Value = -1187.2858
$ws.Range("N31").Value = -2192.8

# CRP row 34: Armoires of the Rich and Famous | Walnut Lumber (item 44023)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 1514
$ws.Range("I34").Value = 1482.2858
$ws.Range("J34").Value = 1602.8
$ws.Range("K34").Value = 1482.2858
$ws.Range("L34").Value = 1602.8
$ws.Range("M34").Value = -1280.2858
$ws.Range("N34").Value = -2006.8

# CRP row 86: Birch, Please | Birch Lumber (item 12584)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H86").Value = 2530.2703
$ws.Range("I86").Value = 1955.2727
$ws.Range("J86").Value = 3373.6
$ws.Range("K86").Value = 1955.2727
$ws.Range("L86").Value = 3373.6
$ws.Range("M86").Value = -832.2727
$ws.Range("N86").Value = -5619.6

# CRP row 89: Built This City on Blocks and Soul (L) | Birch Lumber (item 12584)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H89").Value = 2530.2703
$ws.Range("I89").Value = 1955.2727
$ws.Range("J89").Value = 3373.6
$ws.Range("K89").Value = 9776.363499999999
$ws.Range("L89").Value = 16868
$ws.Range("M89").Value = -4160.363499999999
$ws.Range("N89").Value = -28100

# CUL row 12: Butter Me Up | Kukuru Butter (item 4854)
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 48.9
$ws.Range("I12").Value = 61.46154
$ws.Range("J12").Value = 39.294117
$ws.Range("K12").Value = 184.38462
$ws.Range("L12").Value = 117.882351
$ws.Range("M12").Value = -11.38461999999998
$ws.Range("N12").Value = -463.882351

# CUL row 50: Moving Up in the World | Rolanberry Cheese (item 4725)
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H50").Value = 9230
$ws.Range("I50").Value = 11260
$ws.Range("J50").Value = 95
$ws.Range("K50").Value = 33780
$ws.Range("L50").Value = 285
$ws.Range("M50").Value = -33299
$ws.Range("N50").Value = -1247

# CUL row 53: Rolanberry Fields Forever | Rolanberry Cheese (item 4725)
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H53").Value = 9230
$ws.Range("I53").Value = 11260
$ws.Range("J53").Value = 95
$ws.Range("K53").Value = 33780
$ws.Range("L53").Value = 285
$ws.Range("M53").Value = -33299
$ws.Range("N53").Value = -1247

# CUL row 68: Such a Butter Face | Fermented Butter (item 12895)
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 872.9286
$ws.Range("I68").Value = 837.25
$ws.Range("J68").Value = 920.5
$ws.Range("K68").Value = 2511.75
$ws.Range("L68").Value = 2761.5
$ws.Range("M68").Value = -1700.75
$ws.Range("N68").Value = -4383.5

# CUL row 71: No Margarine of Error (L) | Fermented Butter (item 12895)
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H71").Value = 872.9286
$ws.Range("I71").Value = 837.25
$ws.Range("J71").Value = 920.5
$ws.Range("K71").Value = 7535.25
$ws.Range("L71").Value = 8284.5
$ws.Range("M71").Value = -3479.25
$ws.Range("N71").Value = -16396.5

# CUL row 98: Sweet Kiss of Death | Rice Vinegar (item 19843)
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H98").Value = 612.3077
$ws.Range("I98").Value = 1099.6666
$ws.Range("J98").Value = 466.1
$ws.Range("K98").Value = 3298.9998
$ws.Range("L98").Value = 1398.3
$ws.Range("M98").Value = -1800.9998
$ws.Range("N98").Value = -4394.3

# CUL row 117: A Good Omen | Peppered Popotoes (item 27870)
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H117").Value = 4966
$ws.Range("J117").Value = 5524.875
$ws.Range("L117").Value = 16574.625
$ws.Range("N117").Value = -23458.625

# CUL row 131: The Mountain Steeped | Tsai tou Vounou (item 36060)
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 2130542.5
$ws.Range("I131").Value = 4503.25
$ws.Range("J131").Value = 2859470.2
$ws.Range("K131").Value = 13509.75
$ws.Range("L131").Value = 8578410.600000001
$ws.Range("M131").Value = -8469.75
$ws.Range("N131").Value = -8588490.600000001

# CUL row 132: More Mezcal | Cooking Mezcal (item 43972)
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H132").Value = 35715156
$ws.Range("I132").Value = 52632428
$ws.Range("J132").Value = 919.2222
$ws.Range("K132").Value = 473691852
$ws.Range("L132").Value = 8272.9998
$ws.Range("M132").Value = -473689322
$ws.Range("N132").Value = -13332.9998

# CUL row 133: Friends Are Food | Boiled Alpaca Steak (item 44073)
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H133").Value = 7312.1875
$ws.Range("I133").Value = 3670
$ws.Range("J133").Value = 8152.6924
$ws.Range("K133").Value = 11010
$ws.Range("L133").Value = 24458.0772
$ws.Range("M133").Value = -5950
$ws.Range("N133").Value = -34578.0772

# GSM row 70: Sky Is the Limit | Mythrite Ingot (item 14146)
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 5950
$ws.Range("I70").Value = 5980
$ws.Range("J70").Value = 5900
$ws.Range("K70").Value = 5980
$ws.Range("L70").Value = 5900
$ws.Range("M70").Value = -5710
$ws.Range("N70").Value = -6440

# GSM row 73: Hulls of Broken Dreams (L) | Mythrite Ingot (item 14146)
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H73").Value = 5950
$ws.Range("I73").Value = 5980
$ws.Range("J73").Value = 5900
$ws.Range("K73").Value = 5980
$ws.Range("L73").Value = 5900
$ws.Range("M73").Value = -5044
$ws.Range("N73").Value = -7772

# LTW row 22: Skin off Their Backs | Aldgoat Leather (item 5277)
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 416.66666
$ws.Range("I22").Value = 375
$ws.Range("J22").Value = 500
$ws.Range("K22").Value = 375
$ws.Range("L22").Value = 500
$ws.Range("M22").Value = -80
$ws.Range("N22").Value = -1090

# LTW row 27: Fire and Hide | Aldgoat Leather (item 5277)
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H27").Value = 416.66666
$ws.Range("I27").Value = 375
$ws.Range("J27").Value = 500
$ws.Range("K27").Value = 375
$ws.Range("L27").Value = 500
$ws.Range("M27").Value = -268

# LTW row 40: Best Served Toad | Toad Leather (item 36248)
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 2750
$ws.Range("I40").Value = 2740
$ws.Range("J40").Value = 2775
$ws.Range("K40").Value = 2740
$ws.Range("L40").Value = 2775
$ws.Range("M40").Value = -2604
$ws.Range("N40").Value = -3047
